$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.271.21"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "2.249.18"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.82"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.71%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.74%  "
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0943"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.08%  "
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "2.584.83"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.855"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "2.256.72"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "42.176.10"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +35.79%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0823"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("E36").Value = "  +9.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0317"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "62.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.39%  "
